# Add a new "Croatia" worksheet (Test Data for Croatia Market), cloned from
# the existing "Turkey" sheet and placed right after it as the new last tab.

$wb = $excel.ActiveWorkbook

$turkey = $wb.Worksheets.Item("Turkey")

# Duplicate the Turkey sheet's layout/formatting, inserting the copy right
# after Turkey (i.e. as the new last sheet).
$turkey.Copy($null, $turkey)

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Fill in the market-specific cells. Set B4 before B2 so the new shared
# strings are appended in the expected order.
$croatia.Range("B4").Value = "NGC-3139/T2482"
$croatia.Range("B2").Value = "Croatia Market"

# The Turkey sheet is left with a "select all" selection state after the
# copy/reshuffle.
[void]$turkey.Range("A1:XFD1048576").Select()

# Make Croatia the active sheet/selection, matching the saved workbook state.
# (Done last so Croatia ends up the active/selected tab.)
[void]$croatia.Select()
[void]$croatia.Range("B10").Select()
